# Apply fixes for "un par de cagadas al traer los datos de la J6"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: header rename "Nombre Completo" -> "Jugador"
$ws.Range("A1").Value = "Jugador"

# Row 4 (Franco Saravia): position label Arquero -> Portero
$ws.Range("J4").Value = "Portero"
# J6 - Rendimiento (column V) was blank, should be -2
$ws.Range("V4").Value = -2

# Row 10 (Aldair Fuentes): position label Mediocampista -> Defensa
$ws.Range("J10").Value = "Defensa"

# Row 13 (Jhamir D'Arrigo): fix surname typo and position label
$ws.Range("I13").Value = "D'Arrigo"
$ws.Range("J13").Value = "Defensa"

# Row 20 (Angelo Campos): Arquero -> Portero
$ws.Range("J20").Value = "Portero"

# Row 21 (Angel De la Cruz): Arquero -> Portero
$ws.Range("J21").Value = "Portero"

# Row 26 (Cristian Neira): fill in missing Pos_1 / Pos_2
$ws.Range("C26").Value = "MCO"
$ws.Range("D26").Value = "VLX"
